# Replace the 21 data rows (rows 2-22) of Sheet1 with a new 20-row dataset
# (rows 2-21) reflecting the re-generated / re-sampled "struggle" sensor
# readings, and delete the now-unused last row so the sheet's used range
# shrinks from A1:H22 to A1:H21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, "struggle", 4.229360163211825, -4.647700071334836, -1.350979804992678, 1.217607140541077, -0.0294742472469806, -0.807563841342926),
    @(100, "struggle", 4.624738931655887, -3.006192684173583, -2.288825809955597, 0.5273294448852539, -6.853602409362793, 0.554818332195282),
    @(200, "struggle", -6.357042789459214, -8.575422286987301, 0.2564473152160626, -2.067320823669434, -2.543490171432495, -1.371545195579529),
    @(300, "struggle", 0.5462948679924189, -6.825089752674098, -0.5449948012828849, -1.386664152145386, -2.159103155136108, -0.2104430794715881),
    @(400, "struggle", 0.9978208541870109, -3.706368923187252, -1.51154860854149, -0.4366159439086914, -2.386039733886719, 0.6339253783226013),
    @(500, "struggle", 0.2441467046737658, -3.069634318351747, -2.925750926136973, 0.7061602473258972, 1.022435665130615, 0.9382890462875366),
    @(600, "struggle", 1.612907171249393, -4.853008508682255, -1.383459806442257, -1.432479023933411, 0.08170322328805921, 1.352913856506348),
    @(700, "struggle", -1.189411103725463, -6.66196793317795, 2.149218022823342, -4.308279991149902, 0.7483099102973938, -0.4677700698375702),
    @(800, "struggle", -6.306459784507723, -6.704558491706845, 4.024554014205929, -1.983479499816894, 0.1204931661486625, -1.231657028198242),
    @(900, "struggle", 1.663261890411397, -5.016231019049863, 3.32933139801026, -0.7951938509941101, -1.842216849327088, -2.083966732025146),
    @(1000, "struggle", 6.837078571319559, 0.09169325232505482, 5.010437965393057, -0.8384125232696533, -2.687654256820679, -1.451874017715454),
    @(1100, "struggle", 2.122651159763334, -0.6034613586962226, 2.880795598030089, -0.4280638098716736, -0.8827002644538879, 0.5508477091789246),
    @(1200, "struggle", -0.7601926326751736, 2.327319413423542, 5.726811170578007, 0.1948660165071487, 1.952783465385437, -2.956129550933838),
    @(1300, "struggle", 4.228423535823836, -2.42055988311769, 5.180934607982632, 2.541657686233521, 1.120479583740234, 0.3246748745441437),
    @(1400, "struggle", 2.933720350265484, -5.437817335128782, 5.24878549575806, -3.313027620315552, 2.173916578292847, -4.959309577941895),
    @(1500, "struggle", -2.91136687994004, -3.510188579559316, 5.247701197862618, 0.2906191349029541, -1.299157619476318, 1.514182209968567),
    @(1600, "struggle", -4.280053377151489, -0.6633338928222658, 4.682214915752411, 1.71515691280365, 0.9166033267974854, 0.9940304756164552),
    @(1700, "struggle", -2.581492483615869, 0.5014263689517995, 6.192452192306466, -1.992336988449097, -0.5707008838653564, 0.5294674634933472),
    @(1800, "struggle", -0.986172676086416, 1.025731801986666, -7.288565635681123, -0.0847575515508651, -0.578489363193512, 2.060448408126831),
    @(1900, "struggle", 1.218793094158184, -5.294871598482164, -0.7042694091796555, 1.970193147659302, -0.4350887835025787, 1.578170418739319)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

# The old row 22 is no longer part of the dataset; delete it entirely so
# the sheet's used range becomes A1:H21.
$ws.Rows.Item(22).Delete()
